# Generate Report for Handback
#
# The handback CI run regenerated the report: two source files
# (3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7 / 58c7fef4-7043-4af8-9658-444a7f9a32d6)
# were replaced by new ones (b496a5f4-4a97-4f66-b179-2798b54a02f9 /
# ffff93b24ac0-24d3-45af-9f08-77195abee18a), a new shared commit hash
# (f4805f986f0de95a7dd32215cb1f767a935043db) was produced for both locales'
# handoff/handback xlf files, and new handoff/handback timestamps were
# recorded. Apply this as an exact whole-value substitution across every
# cell and hyperlink display string in the workbook (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# Exact old-value -> new-value pairs (whole cell text replacement, not
# substring), since the two different source GUIDs collapse onto a single
# new xlf filename.
$map = New-Object 'System.Collections.Generic.Dictionary[string,string]'
$map["3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.md"] = "b496a5f4-4a97-4f66-b179-2798b54a02f9.md"
$map["58c7fef4-7043-4af8-9658-444a7f9a32d6.md"] = "ffff93b24ac0-24d3-45af-9f08-77195abee18a.md"
$map["3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.ed74af1e17c37847078fbf243195a30a412ec1b6.zh-cn.xlf"] = "b496a5f4-4a97-4f66-b179-2798b54a02f9.f4805f986f0de95a7dd32215cb1f767a935043db.zh-cn.xlf"
$map["58c7fef4-7043-4af8-9658-444a7f9a32d6.5e2c6408b1e4467cfb04aeec1188a48f0bf3abf9.zh-cn.xlf"] = "b496a5f4-4a97-4f66-b179-2798b54a02f9.f4805f986f0de95a7dd32215cb1f767a935043db.zh-cn.xlf"
$map["3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.ed74af1e17c37847078fbf243195a30a412ec1b6.de-de.xlf"] = "b496a5f4-4a97-4f66-b179-2798b54a02f9.f4805f986f0de95a7dd32215cb1f767a935043db.de-de.xlf"
$map["58c7fef4-7043-4af8-9658-444a7f9a32d6.5e2c6408b1e4467cfb04aeec1188a48f0bf3abf9.de-de.xlf"] = "b496a5f4-4a97-4f66-b179-2798b54a02f9.f4805f986f0de95a7dd32215cb1f767a935043db.de-de.xlf"
$map["2016-03-17 20:49:38"] = "2016-03-17 20:50:56"
$map["2016-03-17 20:49:56"] = "2016-03-17 20:51:31"
$map["2016-03-17 20:49:42"] = "2016-03-17 20:51:03"
$map["2016-03-17 20:50:02"] = "2016-03-17 20:51:36"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string] -and $map.ContainsKey($val)) {
                $cell.Value2 = $map[$val]
            }
        }
    }

    foreach ($hl in $ws.Hyperlinks) {
        $disp = $hl.TextToDisplay
        if ($disp -ne $null -and $map.ContainsKey($disp)) {
            $hl.TextToDisplay = $map[$disp]
        }
    }
}
